# Refresh the cryptos price/volume snapshot (GitHub Actions scheduled update).
# Price (column D) and Volume(1h) (column E) are plain text columns in this
# sheet, so numeric-looking prices are written with a leading apostrophe to
# keep Excel from re-interpreting them as numbers (which would both change
# their cell type and round the displayed text via floating point).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '66.470.24'
$ws.Range('E2').Value = '  +2.81%  '
$ws.Range('D3').Value = '3.188.63'
$ws.Range('E3').Value = '  +1.29%  '
$ws.Range('E4').Value = '  -0.06%  '
$ws.Range('D5').Value = '''596.19'
$ws.Range('E5').Value = '  +3.31%  '
$ws.Range('D6').Value = '''154.61'
$ws.Range('E6').Value = '  +3.65%  '
$ws.Range('E7').Value = '  +0.02%  '
$ws.Range('D8').Value = '''0.564'
$ws.Range('E8').Value = '  +7.26%  '
$ws.Range('D9').Value = '3.183.54'
$ws.Range('E9').Value = '  +1.16%  '
$ws.Range('E10').Value = '  +1.72%  '
$ws.Range('D11').Value = '''5.89'
$ws.Range('E11').Value = '  -3.89%  '
$ws.Range('D12').Value = '''0.519'
$ws.Range('E12').Value = '  +3.61%  '
$ws.Range('D13').Value = '''0.0000269'
$ws.Range('E13').Value = '  +2.46%  '
$ws.Range('D14').Value = '''39.30'
$ws.Range('E14').Value = '  +5.73%  '
$ws.Range('D15').Value = '3.709.52'
$ws.Range('E15').Value = '  +1.30%  '
$ws.Range('E16').Value = '  +4.87%  '
$ws.Range('D17').Value = '66.428.63'
$ws.Range('E17').Value = '  +2.64%  '
$ws.Range('D18').Value = '3.185.19'
$ws.Range('E18').Value = '  +1.35%  '
$ws.Range('E19').Value = '  +0.62%  '
$ws.Range('D20').Value = '''518.44'
$ws.Range('E20').Value = '  +2.73%  '
$ws.Range('D21').Value = '''15.45'
$ws.Range('E21').Value = '  +3.50%  '
$ws.Range('E22').Value = '  +3.59%  '
$ws.Range('D23').Value = '''8.10'
$ws.Range('E23').Value = '  +5.10%  '
$ws.Range('E24').Value = '  -1.88%  '
$ws.Range('D25').Value = '''86.21'
$ws.Range('E25').Value = '  +2.26%  '
$ws.Range('E26').Value = '  +0.09%  '
$ws.Range('D27').Value = '''9.30'
$ws.Range('E27').Value = '  +4.51%  '
$ws.Range('D28').Value = '''3.01'
$ws.Range('E28').Value = '  +3.59%  '
$ws.Range('E29').Value = '  +7.59%  '
$ws.Range('D30').Value = '''7.11'
$ws.Range('E30').Value = '  +14.80%  '
$ws.Range('D31').Value = '''2.93'
$ws.Range('E31').Value = '  +4.08%  '
$ws.Range('D32').Value = '''28.39'
$ws.Range('E32').Value = '  +2.88%  '
$ws.Range('E33').Value = '  +2.97%  '
$ws.Range('E34').Value = '  +0.21%  '
$ws.Range('D35').Value = '''6.53'
$ws.Range('E35').Value = '  +1.17%  '
$ws.Range('D36').Value = '''509.74'
$ws.Range('E36').Value = '  +6.78%  '
$ws.Range('D37').Value = '''54.94'
$ws.Range('E37').Value = '  +0.74%  '
$ws.Range('E38').Value = '  +1.63%  '
$ws.Range('D39').Value = '''0.0427'
$ws.Range('E39').Value = '  +3.10%  '
$ws.Range('E40').Value = '  +9.97%  '
$ws.Range('E41').Value = '  +1.97%  '
$ws.Range('D42').Value = '''2.89'
$ws.Range('E42').Value = '  -0.44%  '
$ws.Range('D43').Value = '''0.303'
$ws.Range('E43').Value = '  +7.50%  '
$ws.Range('D44').Value = '0.0₃0672'
$ws.Range('E44').Value = '  +15.74%  '
$ws.Range('D45').Value = '''2.45'
$ws.Range('E45').Value = '  +1.55%  '
$ws.Range('D46').Value = '2.906.98'
$ws.Range('E46').Value = '  -3.03%  '
$ws.Range('D47').Value = '28.73'
$ws.Range('E47').Value = '  +1.91%  '
$ws.Range('D48').Value = '''0.118'
$ws.Range('E48').Value = '  +3.81%  '
$ws.Range('D49').Value = '''2.69'
$ws.Range('E49').Value = '  +9.14%  '

# Rows 50/51 swapped rank order (ThetaToken now outranks USDe).
$ws.Range('B50').Value = 'ThetaToken'
$ws.Range('C50').Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range('D50').Value = '''2.36'
$ws.Range('E50').Value = '  +5.96%  '
$ws.Range('B51').Value = 'USDe'
$ws.Range('C51').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D51').Value = '0.999'
$ws.Range('E51').Value = '  -0.01%  '
